$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9 (shifts existing rows 9-19 down to 10-20)
$ws.Rows.Item(9).Insert()

# Populate the new row 9 with the new weekly record
$ws.Cells.Item(9, 1).Value = 10
$ws.Cells.Item(9, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(9, 3).Value = "La Araucanía"
$ws.Cells.Item(9, 4).Value = 44427
$ws.Cells.Item(9, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(9, 5).Value = 9
$ws.Cells.Item(9, 6).Value = 100112026
$ws.Cells.Item(9, 7).Value = "Haba"
$ws.Cells.Item(9, 8).Value = "Sin especificar"
$ws.Cells.Item(9, 9).Value = "Primera"
$ws.Cells.Item(9, 10).Value = 30
$ws.Cells.Item(9, 11).Value = 15000
$ws.Cells.Item(9, 12).Value = 15000
$ws.Cells.Item(9, 13).Value = 15000
$ws.Cells.Item(9, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(9, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(9, 16).Value = 600
$ws.Cells.Item(9, 17).Value = 25
$ws.Cells.Item(9, 18).Value = "Hortaliza"
